$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "'30.413.50", "  -0.98%  ")
    ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "'1.920.02", "  +2.84%  ")
    ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "'1.000", "  +0.27%  ")
    ,@(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "'240.82", "  +1.50%  ")
    ,@(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "'1.001", "  +0.23%  ")
    ,@(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "'0.4696", "  -1.18%  ")
    ,@(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "'0.2861", "  +1.51%  ")
    ,@(9, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "'0.06891", "  +6.18%  ")
    ,@(10, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "'107.20", "  +14.79%  ")
    ,@(11, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "'18.36", "  -0.86%  ")
    ,@(12, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "'1.921.60", "  -0.46%  ")
    ,@(13, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "'0.07653", "  +2.29%  ")
    ,@(14, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "'5.189", "  +2.77%  ")
    ,@(15, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "'0.6570", "  +1.83%  ")
    ,@(16, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "'295.25", "  -1.38%  ")
    ,@(17, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "'30.431.98", "  -0.69%  ")
    ,@(18, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "'0.000007655", "  +2.52%  ")
    ,@(19, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "'13.04", "  +0.60%  ")
    ,@(20, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "'1.000", "  +0.12%  ")
    ,@(21, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "'2.157.17", "  +2.57%  ")
    ,@(22, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "'1.001", "  +0.44%  ")
    ,@(23, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "'5.235", "  +1.39%  ")
    ,@(24, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "'6.215", "  +2.20%  ")
    ,@(25, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "'167.81", "  -0.45%  ")
    ,@(26, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "'9.272", "  +0.98%  ")
    ,@(27, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "'21.49", "  +10.56%  ")
    ,@(28, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "'2.045", "  +4.98%  ")
    ,@(29, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "'0.1085", "  +4.24%  ")
    ,@(30, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "'1.360", "  +1.00%  ")
    ,@(31, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "'4.150", "  +1.15%  ")
    ,@(32, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "'3.967", "  +0.58%  ")
    ,@(33, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "'0.05078", "  +2.87%  ")
    ,@(34, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "'0.7432", "  +3.74%  ")
    ,@(35, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "'1.146", "  -2.11%  ")
    ,@(36, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "'2.742", "  +1.39%  ")
    ,@(37, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "'0.02015", "  +4.87%  ")
    ,@(38, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "'2.685", "  -0.82%  ")
    ,@(39, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "'2.059", "  +0.82%  ")
    ,@(40, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "'0.8736", "  -1.55%  ")
    ,@(41, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "'108.44", "  +1.36%  ")
    ,@(42, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "'5.837", "  +5.36%  ")
    ,@(43, "BitcoinSV", "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv", "'53.16", "  +24.44%  ")
    ,@(44, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "'1.000", "  +0.18%  ")
    ,@(45, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "'0.4226", "  +1.05%  ")
    ,@(46, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "'67.73", "  +5.19%  ")
    ,@(47, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "'7.206", "  -2.05%  ")
    ,@(48, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "'9.200", "  +5.44%  ")
    ,@(49, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "'0.1208", "  -1.18%  ")
    ,@(50, "Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "'34.66", "  +0.15%  ")
    ,@(51, "WOONetwork", "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo", "'0.2420", "  +13.73%  ")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
